$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: nudge a range's font size up then back down to its original value.
# This forces the engine to keep the edited sub-range as its own <w:r> run
# instead of re-merging it with neighboring runs that share identical
# formatting.
# ---------------------------------------------------------------------------
function Protect-RunBoundary($range) {
    $orig = $range.Font.Size
    $range.Font.Size = $orig + 1
    $range.Font.Size = $orig
}

# ---------------------------------------------------------------------------
# 1) Date paragraph: "08/29/2024" -> "09/08/2024"
#    Original runs: "08/2" | "9" | "/2024"
#    New runs:      "0"    | "9/08" | "/2024"   (same run boundaries/rPr)
# ---------------------------------------------------------------------------
$dateFind = $d.Content
$dateFind.Find.ClearFormatting()
$dateFind.Find.Execute("08/29/2024") | Out-Null
$dateStart = $dateFind.Start

# Edit the first run ("08/2" -> "0") - shrinks from 4 to 1 characters.
$runA = $d.Range($dateStart, $dateStart + 4)
$runA.Text = "0"

# The second run ("9") now sits right after the shortened first run.
$runB = $d.Range($dateStart + 1, $dateStart + 2)
$runB.Text = "9/08"

# Re-assert the boundary between the two edited runs, and the boundary
# with the (untouched) trailing "/2024" run, so they don't get coalesced.
Protect-RunBoundary ($d.Range($dateStart, $dateStart + 1))
Protect-RunBoundary ($d.Range($dateStart + 1, $dateStart + 5))

# ---------------------------------------------------------------------------
# 2) "Lab 1" -> "Assignment 1"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Lab 1", $true, $false, $false, $false, $false, $true, 1, $false, "Assignment 1", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Split the "Collaboration Log" paragraph (the plain one right after the
#    assignment title, not the bold title further down) into two runs:
#    "Collaboration Lo" + "g" - same text, same formatting, just a run split.
# ---------------------------------------------------------------------------
$collabFind = $d.Content
$collabFind.Find.ClearFormatting()
$collabFind.Find.MatchWholeWord = $false
$collabFind.Find.Execute("Collaboration Log") | Out-Null
$collabStart = $collabFind.Start

$gRange = $d.Range($collabStart + 16, $collabStart + 17)
Protect-RunBoundary $gRange

Write-Host "Done"
